# Update column G ("K") values to reflect the regenerated save_data
# (switch from Strike# to K, with recalculated std/mean-derived s_vals)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 0
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 2
